# Auto-generated Excel COM-interop script applying diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F23").Value = 38
$ws.Range("G23").Value = 1557.62
$ws.Range("F24").Value = 35
$ws.Range("G24").Value = 3585.75
$ws.Range("F26").Value = 59
$ws.Range("G26").Value = 1510.99
$ws.Range("F27").Value = 81
$ws.Range("G27").Value = 2904.66
$ws.Range("F28").Value = 60
$ws.Range("G28").Value = 1844.4
$ws.Range("B34").Value = 68408.46000000001
$ws.Range("F38").Value = 517
$ws.Range("G38").Value = 18823.97
$ws.Range("F42").Value = 84
$ws.Range("G42").Value = 3538.08
$ws.Range("F44").Value = 41
$ws.Range("G44").Value = 1448.12
$ws.Range("F46").Value = 65
$ws.Range("G46").Value = 2416.7
$ws.Range("F55").Value = 144
$ws.Range("G55").Value = 8029.44
$ws.Range("F64").Value = 72
$ws.Range("G64").Value = 5725.44
$ws.Range("B66").Value = 248343.23
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("B83").Value = 47685.21
$ws.Range("F95").Value = 6
$ws.Range("G95").Value = 2230.8
$ws.Range("B97").Value = 20437.89
$ws.Range("F111").Value = 15
$ws.Range("G111").Value = 955.5
$ws.Range("F114").Value = 97
$ws.Range("G114").Value = 4533.78
$ws.Range("F116").Value = 45
$ws.Range("G116").Value = 2797.65
$ws.Range("B123").Value = 85075.52
$ws.Range("F132").Value = 1
$ws.Range("G132").Value = 1690.08
$ws.Range("B133").Value = 9565.219999999999
$ws.Range("F141").Value = 76
$ws.Range("G141").Value = 3956.56
$ws.Range("B147").Value = 26352.67
$ws.Range("F172").Value = 124
$ws.Range("G172").Value = 7876.48
$ws.Range("F173").Value = 70
$ws.Range("G173").Value = 5502.7
$ws.Range("F174").Value = 44
$ws.Range("G174").Value = 7055.84
$ws.Range("F178").Value = 111
$ws.Range("G178").Value = 10740.36
$ws.Range("B193").Value = 81612.16
$ws.Range("F210").Value = 152
$ws.Range("G210").Value = 8270.32
$ws.Range("B218").Value = 94605.92999999999
$ws.Range("F220").Value = 65
$ws.Range("G220").Value = 4106.7
$ws.Range("F222").Value = 1609
$ws.Range("G222").Value = 29766.5
$ws.Range("F227").Value = 75
$ws.Range("G227").Value = 8595
$ws.Range("B229").Value = 46703.11
$ws.Range("F263").Value = 24
$ws.Range("G263").Value = 2488.8
$ws.Range("F264").Value = 120
$ws.Range("G264").Value = 4180.8
$ws.Range("F275").Value = 6
$ws.Range("G275").Value = 455.46
$ws.Range("F278").Value = 57
$ws.Range("G278").Value = 7725.78
$ws.Range("B290").Value = 66194
$ws.Range("C290").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F290").Value = 35
$ws.Range("G290").Value = 2998.8
$ws.Range("B291").Value = 64983
$ws.Range("C291").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F291").Value = 6
$ws.Range("G291").Value = 514.08
$ws.Range("F293").Value = 25
$ws.Range("G293").Value = 2192.5
$ws.Range("B295").Value = 158774.62
$ws.Range("F307").Value = 34
$ws.Range("G307").Value = 4878.32
$ws.Range("B328").Value = 23534.51
$ws.Range("F368").Value = 68
$ws.Range("G368").Value = 2177.36
$ws.Range("B372").Value = 77516.28999999999
$ws.Range("F387").Value = 514
$ws.Range("G387").Value = 49652.4
$ws.Range("B389").Value = 67421.17
$ws.Range("F408").Value = 259
$ws.Range("G408").Value = 4105.15
$ws.Range("B417").Value = 192551.39
$ws.Range("F452").Value = 64
$ws.Range("G452").Value = 17358.08
$ws.Range("B458").Value = 120034.75
$ws.Range("F466").Value = 37
$ws.Range("G466").Value = 1524.4
$ws.Range("B476").Value = 57127.36
$ws.Range("B496").Value = 60025
$ws.Range("E496").Value = 37.22
$ws.Range("F496").Value = -98
$ws.Range("G496").Value = -3217.34
$ws.Range("B497").Value = 64833
$ws.Range("E497").Value = 34.9
$ws.Range("F497").Value = 88
$ws.Range("G497").Value = 2889.04
$ws.Range("B506").Value = 64830
$ws.Range("E506").Value = 34.9
$ws.Range("F506").Value = 88
$ws.Range("G506").Value = 2889.04
$ws.Range("B507").Value = 60022
$ws.Range("E507").Value = 37.22
$ws.Range("F507").Value = -113
$ws.Range("G507").Value = -3709.79
$ws.Range("F512").Value = 44
$ws.Range("G512").Value = 5217.52
$ws.Range("F519").Value = 484
$ws.Range("G519").Value = 26561.92
$ws.Range("F520").Value = 76
$ws.Range("G520").Value = 2082.4
$ws.Range("F523").Value = 183
$ws.Range("G523").Value = 15666.63
$ws.Range("B525").Value = 150507.26
$ws.Range("F528").Value = 343
$ws.Range("G528").Value = 5439.98
$ws.Range("B535").Value = 33310.54
$ws.Range("F558").Value = 274
$ws.Range("G558").Value = 33386.9
$ws.Range("F560").Value = 58
$ws.Range("G560").Value = 4662.04
$ws.Range("B561").Value = 40439.56
$ws.Range("F567").Value = 2
$ws.Range("G567").Value = 299.76
$ws.Range("B573").Value = 41652.84
$ws.Range("F605").Value = 220
$ws.Range("G605").Value = 29282
$ws.Range("B607").Value = 29822.04
$ws.Range("F609").Value = 72
$ws.Range("G609").Value = 7834.32
$ws.Range("F616").Value = 33
$ws.Range("G616").Value = 4711.41
$ws.Range("F625").Value = 372
$ws.Range("G625").Value = 13700.76
$ws.Range("B628").Value = 255369.74
$ws.Range("F673").Value = 52
$ws.Range("G673").Value = 1571.96
$ws.Range("F674").Value = 1345
$ws.Range("G674").Value = 219382.95
$ws.Range("B680").Value = 236701.19
$ws.Range("F703").Value = 11
$ws.Range("G703").Value = 5337.97
$ws.Range("B713").Value = 82589.78999999999
$ws.Range("B718").Value = 3620229.31
$ws.Range("B719").Value = 3620229.31
